$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.336.36"
Set-TextValue $ws.Range("E2") "  -3.41%  "
Set-TextValue $ws.Range("D3") "1.933.22"
Set-TextValue $ws.Range("E3") "  -3.84%  "
Set-TextValue $ws.Range("E4") "  +0.09%  "
Set-TextValue $ws.Range("D5") "249.02"
Set-TextValue $ws.Range("E5") "  -3.62%  "
Set-TextValue $ws.Range("D6") "0.7121"
Set-TextValue $ws.Range("E6") "  -6.64%  "
Set-TextValue $ws.Range("D7") "1.000"
Set-TextValue $ws.Range("E7") "  +0.11%  "
Set-TextValue $ws.Range("D8") "0.3280"
Set-TextValue $ws.Range("E8") "  -8.92%  "
Set-TextValue $ws.Range("D9") "27.64"
Set-TextValue $ws.Range("E9") "  -4.68%  "
Set-TextValue $ws.Range("D10") "0.06854"
Set-TextValue $ws.Range("E10") "  -3.31%  "
Set-TextValue $ws.Range("D11") "0.8070"
Set-TextValue $ws.Range("E11") "  -5.67%  "
Set-TextValue $ws.Range("D12") "0.08121"
Set-TextValue $ws.Range("E12") "  +0.00%  "
Set-TextValue $ws.Range("D13") "1.932.56"
Set-TextValue $ws.Range("E13") "  -3.83%  "
Set-TextValue $ws.Range("D14") "5.446"
Set-TextValue $ws.Range("E14") "  -3.16%  "
Set-TextValue $ws.Range("D15") "95.05"
Set-TextValue $ws.Range("E15") "  -7.03%  "
Set-TextValue $ws.Range("D16") "14.61"
Set-TextValue $ws.Range("E16") "  +0.20%  "
Set-TextValue $ws.Range("D17") "259.90"
Set-TextValue $ws.Range("E17") "  -5.19%  "
Set-TextValue $ws.Range("D18") "30.326.44"
Set-TextValue $ws.Range("E18") "  -3.42%  "
Set-TextValue $ws.Range("D19") "0.000008046"
Set-TextValue $ws.Range("E19") "  +0.39%  "
Set-TextValue $ws.Range("D20") "5.848"
Set-TextValue $ws.Range("E20") "  -1.69%  "
Set-TextValue $ws.Range("D21") "2.186.72"
Set-TextValue $ws.Range("E21") "  -3.45%  "
Set-TextValue $ws.Range("D22") "1.000"
Set-TextValue $ws.Range("E22") "  +0.10%  "
Set-TextValue $ws.Range("D23") "0.9995"
Set-TextValue $ws.Range("D24") "6.882"
Set-TextValue $ws.Range("E24") "  -6.06%  "
Set-TextValue $ws.Range("D25") "9.727"
Set-TextValue $ws.Range("E25") "  -4.22%  "
Set-TextValue $ws.Range("D26") "159.73"
Set-TextValue $ws.Range("E26") "  -2.58%  "
Set-TextValue $ws.Range("D27") "2.371"
Set-TextValue $ws.Range("E27") "  +0.47%  "
Set-TextValue $ws.Range("D28") "19.15"
Set-TextValue $ws.Range("E28") "  -4.78%  "
Set-TextValue $ws.Range("D29") "0.1333"
Set-TextValue $ws.Range("E29") "  -8.81%  "
Set-TextValue $ws.Range("D30") "1.562"
Set-TextValue $ws.Range("E30") "  -3.37%  "
Set-TextValue $ws.Range("E31") "  -0.26%  "
Set-TextValue $ws.Range("D32") "4.418"
Set-TextValue $ws.Range("E32") "  -5.02%  "
Set-TextValue $ws.Range("D33") "4.223"
Set-TextValue $ws.Range("E33") "  -4.31%  "
Set-TextValue $ws.Range("D34") "0.05095"
Set-TextValue $ws.Range("E34") "  -2.42%  "
Set-TextValue $ws.Range("D35") "1.224"
Set-TextValue $ws.Range("E35") "  -1.13%  "
Set-TextValue $ws.Range("D36") "0.7430"
Set-TextValue $ws.Range("E36") "  -2.88%  "
Set-TextValue $ws.Range("D37") "2.769"
Set-TextValue $ws.Range("E37") "  -0.99%  "
Set-TextValue $ws.Range("D38") "0.01982"
Set-TextValue $ws.Range("E38") "  -2.09%  "
Set-TextValue $ws.Range("D39") "2.819"
Set-TextValue $ws.Range("E39") "  -4.32%  "
Set-TextValue $ws.Range("D40") "79.47"
Set-TextValue $ws.Range("E40") "  -1.86%  "
Set-TextValue $ws.Range("D41") "6.601"
Set-TextValue $ws.Range("E41") "  -2.54%  "
Set-TextValue $ws.Range("D42") "0.4475"
Set-TextValue $ws.Range("E42") "  -6.47%  "
Set-TextValue $ws.Range("D43") "2.001"
Set-TextValue $ws.Range("E43") "  -8.86%  "
Set-TextValue $ws.Range("D44") "1.001"
Set-TextValue $ws.Range("E44") "  +0.10%  "
Set-TextValue $ws.Range("D45") "0.8366"
Set-TextValue $ws.Range("E45") "  -3.12%  "
Set-TextValue $ws.Range("D46") "102.06"
Set-TextValue $ws.Range("E46") "  -2.82%  "
Set-TextValue $ws.Range("D47") "9.779"
Set-TextValue $ws.Range("E47") "  -2.56%  "
Set-TextValue $ws.Range("D48") "7.329"
Set-TextValue $ws.Range("E48") "  -4.93%  "
Set-TextValue $ws.Range("D49") "36.56"
Set-TextValue $ws.Range("E49") "  -1.24%  "
Set-TextValue $ws.Range("D50") "1.482"
Set-TextValue $ws.Range("E50") "  -0.23%  "

# Row 51: Decentraland -> Cronos
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.05949"
Set-TextValue $ws.Range("E51") "  -0.51%  "
